# Applies the "Added Keras TF files" edit:
#  - Basic Learning ("FULL SET RF min-max" block, rows 69-73): reorders / updates
#    the RF-* results and adds a new "RF-250-entropy" result row.
#  - Advanced Learning ("Numeric Only min-max 1" block, rows 37-40): labels the
#    previously-blank header and fills in the first three result rows, plus a
#    couple of small value corrections a few rows above it.
#  - Selection / active sheet bookkeeping to match the saved UI state.

$wb = $excel.ActiveWorkbook

$wsBasic = $wb.Worksheets.Item("Basic Learning")
$wsAdv   = $wb.Worksheets.Item("Advanced Learning")

# ---------------------------------------------------------------------------
# Basic Learning: "FULL SET RF min-max" table (rows 69-73)
# ---------------------------------------------------------------------------
$wsBasic.Range("A69").Value = "RF-250"
$wsBasic.Range("B69").Value = 0.94267617443850704
$wsBasic.Range("C69").Value = 0.94
$wsBasic.Range("D69").Value = 0.94
$wsBasic.Range("E69").Value = 0.94

$wsBasic.Range("A70").Value = "RF-500"
$wsBasic.Range("B70").Value = 0.94284942100995806
$wsBasic.Range("C70").Value = 0.94
$wsBasic.Range("D70").Value = 0.94
$wsBasic.Range("E70").Value = 0.94

$wsBasic.Range("A71").Value = "RF-1000"
$wsBasic.Range("B71").Value = 0.94316819470142599
$wsBasic.Range("C71").Value = 0.94
$wsBasic.Range("D71").Value = 0.94
$wsBasic.Range("E71").Value = 0.94

$wsBasic.Range("A72").Value = "RF-100-entropy"
$wsBasic.Range("B72").Value = 0.94522636397025706
$wsBasic.Range("C72").Value = 0.95
$wsBasic.Range("D72").Value = 0.95
$wsBasic.Range("E72").Value = 0.94

$wsBasic.Range("A73").Value = "RF-250-entropy"
$wsBasic.Range("B73").Value = 0.950680166039514
$wsBasic.Range("C73").Value = 0.95
$wsBasic.Range("D73").Value = 0.95
$wsBasic.Range("E73").Value = 0.95

# ---------------------------------------------------------------------------
# Advanced Learning: couple of value corrections (rows 22-24)
# ---------------------------------------------------------------------------
$wsAdv.Range("B22").Value = 0.64724470562146497
$wsAdv.Range("D22").Value = 0.65

$wsAdv.Range("B23").Value = 0.65204706730236095
$wsAdv.Range("E23").Value = 0.63

$wsAdv.Range("B24").Value = 0.65145803304135697

# ---------------------------------------------------------------------------
# Advanced Learning: "Numeric Only min-max 1" table (rows 37-40)
# ---------------------------------------------------------------------------
$wsAdv.Range("A37").Value = "Numeric Only min-max 1"

$wsAdv.Range("B38").Value = 0.66467318993236502
$wsAdv.Range("C38").Value = 0.64
$wsAdv.Range("D38").Value = 0.65
$wsAdv.Range("E38").Value = 0.63

$wsAdv.Range("B39").Value = 0.803803082381638
$wsAdv.Range("C39").Value = 0.8
$wsAdv.Range("D39").Value = 0.8
$wsAdv.Range("E39").Value = 0.8

$wsAdv.Range("B40").Value = 0.74379088590752795
$wsAdv.Range("C40").Value = 0.74
$wsAdv.Range("D40").Value = 0.74
$wsAdv.Range("E40").Value = 0.74

# Column A on Advanced Learning widens (auto-fit) once it holds the longer
# "Numeric Only min-max 1" label.
$wsAdv.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------------
# Recalculate so every dependent formula (Basic/Advanced Learning summary
# columns, ALL Learning sheet, etc.) picks up the new cached values.
# ---------------------------------------------------------------------------
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------------
# UI state: the saved workbook shows "Advanced Learning" as the active sheet
# with a fresh selection/scroll position, and "Basic Learning" scrolled to a
# new top-left cell with its own selection.
# ---------------------------------------------------------------------------
$wsBasic.Range("A55").Select()
$wsBasic.Range("B73").Select()

$wsAdv.Activate()
$wsAdv.Range("A16").Select()
$wsAdv.Range("E40").Select()
